# Automatische test-sync: 2025-07-22 17:26:50
# Adds the new test-mail row (row 20) to the "Logs" sheet, grows the
# conditional-formatting ranges to include it, and refreshes the
# "Dashboard" category-count table to reflect the new "Openingstijden /
# Locatie" entry.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new log row (row 20) -----------------------------------
$logs.Cells.Item(20, 1).Value  = "Wat zijn jullie openingstijden?"
$logs.Cells.Item(20, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item(20, 3).Value  = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Cells.Item(20, 4).Value  = "Openingstijden / Locatie"
$logs.Cells.Item(20, 5).Value  = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item(20, 6).Value  = "2025-07-22 17:26:32"
$logs.Cells.Item(20, 7).Value  = "Ja"
$logs.Cells.Item(20, 8).Value  = "Nee"
$logs.Cells.Item(20, 9).Value  = "Ja"
$logs.Cells.Item(20, 10).Value = "Ja"

# The multi-line text in E20 makes the runtime auto-mark the row with an
# explicit (custom) height; re-auto-fitting brings it back to the sheet's
# default row height and drops the custom-height flag again, matching the
# rest of the sheet (no row in the original file carries an explicit
# height).
$logs.Rows.Item(20).AutoFit()

# --- 2. Grow the conditional-formatting ranges so row 20 is included -------
$rD = $logs.Range("D2:D20")
$rD.FormatConditions.Item(1).ModifyAppliesToRange($rD)

$rG = $logs.Range("G2:G20")
$rG.FormatConditions.Item(1).ModifyAppliesToRange($rG)

$rH = $logs.Range("H2:H20")
$rH.FormatConditions.Item(1).ModifyAppliesToRange($rH)

$rI = $logs.Range("I2:I20")
$rI.FormatConditions.Item(1).ModifyAppliesToRange($rI)

$rJ = $logs.Range("J2:J20")
$rJ.FormatConditions.Item(1).ModifyAppliesToRange($rJ)

# --- 3. Refresh the Dashboard category counts -------------------------------
# The new mail belongs to "Openingstijden / Locatie", bumping its count
# from 1 to 2, which ties it with "Intern verzoek / Actie voor medewerker".
# The rows are reordered (descending by count) accordingly:
#   row 4: Openingstijden / Locatie        -> 2
#   row 5: Overig                          -> 2 (unchanged)
#   row 6: Intern verzoek / Actie medewerker -> 2 (was 1 occurrence spot, now count 2)
$dash.Cells.Item(4, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(4, 2).Value = 2

$dash.Cells.Item(6, 1).Value = "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(6, 2).Value = 2
